$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D4").Value = 96.95999200548977
$ws.Range("E4").Value = 51.43841178994626
$ws.Range("F4").Value = 0.5305117164926527
$ws.Range("G4").Value = 1.884972506566402
$ws.Range("H4").Value = 141.1197404945415
$ws.Range("I4").Value = 0.02437799899780657
$ws.Range("J4").Value = 0.3097209488469161
$ws.Range("K4").Value = 4.849417570978403
$ws.Range("L4").Value = 1.436524677090347
$ws.Range("M4").Value = 5.623218013439327
$ws.Range("N4").Value = 0.2179608714068308
$ws.Range("O4").Value = 0.1999594387507386
$ws.Range("P4").Value = 0.0008709411613381235
$ws.Range("Q4").Value = 28.30095673212782
$ws.Range("R4").Value = 5.496810862794518
$ws.Range("S4").Value = 51.23543111188337
$ws.Range("T4").Value = -1410.979444074008
$ws.Range("U4").Value = -0.04382055122732709
$ws.Range("V4").Value = -3.096338547307823
$ws.Range("W4").Value = -20.19321897765622
$ws.Range("X4").Value = -50.73536927159876
$ws.Range("Y4").Value = 36.8701843409799

$ws.Range("D5").Value = 96.95999200548977
$ws.Range("E5").Value = 51.43158450908959
$ws.Range("F5").Value = 0.530441303111675
$ws.Range("G5").Value = 1.885222727064803
$ws.Range("H5").Value = 141.1197174044064
$ws.Range("I5").Value = 0.02436912961275084
$ws.Range("J5").Value = 0.3097208946055616
$ws.Range("K5").Value = 4.844723212067038
$ws.Range("L5").Value = 1.43614204833284
$ws.Range("M5").Value = 5.623095151036978
$ws.Range("N5").Value = 0.2179608714068308
$ws.Range("O5").Value = 0.1999594387507386
$ws.Range("P5").Value = 0.0008709411613381235
$ws.Range("Q5").Value = 28.30095673212782
$ws.Range("R5").Value = 5.496810862794518
$ws.Range("S5").Value = 51.23543111188337
$ws.Range("T5").Value = -1410.979213172657
$ws.Range("U5").Value = -0.04373185737676977
$ws.Range("V5").Value = -3.096338004894278
$ws.Range("W5").Value = -20.14627538854256
$ws.Range("X5").Value = -50.73414064757526
$ws.Range("Y5").Value = 36.87401062855497

$ws.Range("D6").Value = 96.95999200548977
$ws.Range("E6").Value = 51.42475722637028
$ws.Range("F6").Value = 0.530370889711487
$ws.Range("G6").Value = 1.88547301407131
$ws.Range("H6").Value = 141.1196943142422
$ws.Range("I6").Value = 0.02436026022587612
$ws.Range("J6").Value = 0.3097208403623881
$ws.Range("K6").Value = 4.840028852922842
$ws.Range("L6").Value = 1.435759419109672
$ws.Range("M6").Value = 5.622972288634628
$ws.Range("N6").Value = 0.2179608714068308
$ws.Range("O6").Value = 0.1999594387507386
$ws.Range("P6").Value = 0.0008709411613381235
$ws.Range("Q6").Value = 28.30095673212782
$ws.Range("R6").Value = 5.496810862794518
$ws.Range("S6").Value = 51.23543111188337
$ws.Range("T6").Value = -1410.978982271015
$ws.Range("U6").Value = -0.04364316350802255
$ws.Range("V6").Value = -3.096337462462543
$ws.Range("W6").Value = -20.0993317971006
$ws.Range("X6").Value = -50.73291202355176
$ws.Range("Y6").Value = 36.87783692078665

$ws.Range("D7").Value = 96.95999200548977
$ws.Range("E7").Value = 51.41792994551361
$ws.Range("F7").Value = 0.5303004763305094
$ws.Range("G7").Value = 1.885723367475821
$ws.Range("H7").Value = 141.119671224078
$ws.Range("I7").Value = 0.02435139084082039
$ws.Range("J7").Value = 0.3097207861219431
$ws.Range("K7").Value = 4.835334494244307
$ws.Range("L7").Value = 1.435376790584996
$ws.Range("M7").Value = 5.622849426232278
$ws.Range("N7").Value = 0.2179608714068308
$ws.Range("O7").Value = 0.1999594387507386
$ws.Range("P7").Value = 0.0008709411613381235
$ws.Range("Q7").Value = 28.30095673212782
$ws.Range("R7").Value = 5.496810862794518
$ws.Range("S7").Value = 51.23543111188337
$ws.Range("T7").Value = -1410.978751369374
$ws.Range("U7").Value = -0.04355446965746523
$ws.Range("V7").Value = -3.096336920058093
$ws.Range("W7").Value = -20.05238821031526
$ws.Range("X7").Value = -50.73168339952826
$ws.Range("Y7").Value = 36.88166320603341

$ws.Range("D8").Value = 96.95999200548977
$ws.Range("E8").Value = 51.41110266465694
$ws.Range("F8").Value = 0.5302300629495318
$ws.Range("G8").Value = 1.885973787373089
$ws.Range("H8").Value = 141.119648133943
$ws.Range("I8").Value = 0.02434252145576465
$ws.Range("J8").Value = 0.3097207318796791
$ws.Range("K8").Value = 4.830640134867281
$ws.Range("L8").Value = 1.43499416206032
$ws.Range("M8").Value = 5.622726563364267
$ws.Range("N8").Value = 0.2179608714068308
$ws.Range("O8").Value = 0.1999594387507386
$ws.Range("P8").Value = 0.0008709411613381235
$ws.Range("Q8").Value = 28.30095673212782
$ws.Range("R8").Value = 5.496810862794518
$ws.Range("S8").Value = 51.23543111188337
$ws.Range("T8").Value = -1410.978520468023
$ws.Range("U8").Value = -0.04346577580690791
$ws.Range("V8").Value = -3.096336377635453
$ws.Range("W8").Value = -20.00544461654499
$ws.Range("X8").Value = -50.73045477084816
$ws.Range("Y8").Value = 36.88548949128017

$ws.Range("D9").Value = 96.95999200548977
$ws.Range("E9").Value = 51.40427538286895
$ws.Range("F9").Value = 0.5301596495589489
$ws.Range("G9").Value = 1.886224273823784
$ws.Range("H9").Value = 141.1196250437642
$ws.Range("I9").Value = 0.02433365207616589
$ws.Range("J9").Value = 0.3097206776392341
$ws.Range("K9").Value = 4.825945775955915
$ws.Range("L9").Value = 1.434611532604322
$ws.Range("M9").Value = 5.622603700961918
$ws.Range("N9").Value = 0.2179608714068308
$ws.Range("O9").Value = 0.1999594387507386
$ws.Range("P9").Value = 0.0008709411613381235
$ws.Range("Q9").Value = 28.30095673212782
$ws.Range("R9").Value = 5.496810862794518
$ws.Range("S9").Value = 51.23543111188337
$ws.Range("T9").Value = -1410.978289566236
$ws.Range("U9").Value = -0.04337708201092028
$ws.Range("V9").Value = -3.096335835231002
$ws.Range("W9").Value = -19.95850102743134
$ws.Range("X9").Value = -50.72922614682466
$ws.Range("Y9").Value = 36.88931578584015
